# Einheitliche Namengebung in Dokumentation
#
# 1) Fix a wrong date in row 5 (Kontrollschaltung.pdf: 2017-10-24 -> 2017-11-09)
# 2) Split the combined filename in row 7 ("Alles Da.jpg; Angekommen.jpeg")
#    into two separate rows: row 7 keeps "Angekommen.jpeg", and a new row is
#    inserted right after it for "Alles Da.jpg;" with its own description.
# 3) Insert a new row for a previously undocumented milestone
#    ("2017-12.11.jpg" / Raspberry Pi Gehäuse) after the i2c_protokoll.jpg row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) correct the date on row 5 -----------------------------------------
$ws.Range("A5").Value = 43048

# --- 2) split row 7 into two rows ------------------------------------------
$ws.Range("B7").Value = "Angekommen.jpeg"

$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = 43058
$ws.Range("B8").Value = "Alles Da.jpg;"
$ws.Range("C8").Value = "Die Excel Tabelle"

# --- 3) add the missing Raspberry-Pi-Gehäuse milestone ---------------------
# After step 2 the "i2c_protokoll.jpg" entry (old row 17) now lives on row 18.
$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = 43080
$ws.Range("B19").Value = "2017-12.11.jpg"
$ws.Range("C19").Value = "Raspberry Pi gehäuse fertig. Testdruck geschaftt. Etwas verspätet!"

# --- cosmetic: restore the scroll position / selection seen in the file ----
$ws.Range("B39").Select()
